$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddOpportunity")

# Remove the trailing formatting-only rows (5-9) so the used range shrinks to A1:AI4
$ws.Range("A5:AI9").EntireRow.Delete()

# Rename "Tec Alliance" -> "TEC Canada" and turn it into a hyperlink in both data rows
$wb.Hyperlinks.Add($ws.Range("AG2"), "https://hl--test.sandbox.my.salesforce.com/0015A00002LZtzTQAT", "", "", "TEC Canada")
$wb.Hyperlinks.Add($ws.Range("AG3"), "https://hl--test.sandbox.my.salesforce.com/0015A00002LZtzTQAT", "", "", "TEC Canada")

# Leave the selection where the author ended up working
$null = $ws.Range("AG3").Select()
